{"js": "// Update the 25 two-digit multiplication problems/answers in the table\n// from their original values to the newly generated ones.\nconst replacements = [\n  [\"23\u00d767=1541\", \"72\u00d770=5040\"],\n  [\"41\u00d786=3526\", \"85\u00d744=3740\"],\n  [\"11\u00d797=1067\", \"49\u00d721=1029\"],\n  [\"40\u00d774=2960\", \"49\u00d730=1470\"],\n  [\"97\u00d795=9215\", \"32\u00d769=2208\"],\n  [\"31\u00d717=527\", \"41\u00d734=1394\"],\n  [\"56\u00d784=4704\", \"53\u00d713=689\"],\n  [\"22\u00d787=1914\", \"49\u00d737=1813\"],\n  [\"47\u00d748=2256\", \"12\u00d743=516\"],\n  [\"22\u00d768=1496\", \"94\u00d742=3948\"],\n  [\"12\u00d749=588\", \"64\u00d754=3456\"],\n  [\"23\u00d725=575\", \"84\u00d730=2520\"],\n  [\"77\u00d722=1694\", \"35\u00d754=1890\"],\n  [\"51\u00d715=765\", \"42\u00d741=1722\"],\n  [\"18\u00d753=954\", \"97\u00d749=4753\"],\n  [\"94\u00d733=3102\", \"47\u00d755=2585\"],\n  [\"14\u00d761=854\", \"40\u00d772=2880\"],\n  [\"87\u00d728=2436\", \"76\u00d751=3876\"],\n  [\"39\u00d789=3471\", \"39\u00d785=3315\"],\n  [\"76\u00d791=6916\", \"98\u00d780=7840\"],\n  [\"40\u00d798=3920\", \"49\u00d731=1519\"],\n  [\"65\u00d724=1560\", \"17\u00d711=187\"],\n  [\"41\u00d797=3977\", \"43\u00d770=3010\"],\n  [\"28\u00d745=1260\", \"59\u00d733=1947\"],\n  [\"21\u00d732=672\", \"41\u00d740=1640\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 two-digit multiplication problems/answers in the table\n# from their original values to the newly generated ones.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"23\u00d767=1541\", \"72\u00d770=5040\"),\n    @(\"41\u00d786=3526\", \"85\u00d744=3740\"),\n    @(\"11\u00d797=1067\", \"49\u00d721=1029\"),\n    @(\"40\u00d774=2960\", \"49\u00d730=1470\"),\n    @(\"97\u00d795=9215\", \"32\u00d769=2208\"),\n    @(\"31\u00d717=527\", \"41\u00d734=1394\"),\n    @(\"56\u00d784=4704\", \"53\u00d713=689\"),\n    @(\"22\u00d787=1914\", \"49\u00d737=1813\"),\n    @(\"47\u00d748=2256\", \"12\u00d743=516\"),\n    @(\"22\u00d768=1496\", \"94\u00d742=3948\"),\n    @(\"12\u00d749=588\", \"64\u00d754=3456\"),\n    @(\"23\u00d725=575\", \"84\u00d730=2520\"),\n    @(\"77\u00d722=1694\", \"35\u00d754=1890\"),\n    @(\"51\u00d715=765\", \"42\u00d741=1722\"),\n    @(\"18\u00d753=954\", \"97\u00d749=4753\"),\n    @(\"94\u00d733=3102\", \"47\u00d755=2585\"),\n    @(\"14\u00d761=854\", \"40\u00d772=2880\"),\n    @(\"87\u00d728=2436\", \"76\u00d751=3876\"),\n    @(\"39\u00d789=3471\", \"39\u00d785=3315\"),\n    @(\"76\u00d791=6916\", \"98\u00d780=7840\"),\n    @(\"40\u00d798=3920\", \"49\u00d731=1519\"),\n    @(\"65\u00d724=1560\", \"17\u00d711=187\"),\n    @(\"41\u00d797=3977\", \"43\u00d770=3010\"),\n    @(\"28\u00d745=1260\", \"59\u00d733=1947\"),\n    @(\"21\u00d732=672\", \"41\u00d740=1640\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
